$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 123.85714
$ws.Range("I5").Value = 123.85714
$ws.Range("K5").Value = 123.85714
$ws.Range("M5").Value = -8.857140000000001

$ws.Range("H6").Value = 200073.4
$ws.Range("I6").Value = 200073.4
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 600220.2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -600108.2
$ws.Range("N6").ClearContents()

$ws.Range("H12").Value = 164.5
$ws.Range("I12").Value = 173.71428
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 173.71428
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -3.714280000000002
$ws.Range("N12").Value = -440

$ws.Range("H40").Value = 5080.5884
$ws.Range("I40").Value = 2960
$ws.Range("J40").Value = 5363.3335
$ws.Range("K40").Value = 2960
$ws.Range("L40").Value = 5363.3335
$ws.Range("M40").Value = -2785
$ws.Range("N40").Value = -5713.3335

$ws.Range("H64").Value = 6699
$ws.Range("I64").Value = 4500
$ws.Range("J64").Value = 9997.5
$ws.Range("K64").Value = 4500
$ws.Range("L64").Value = 9997.5
$ws.Range("M64").Value = -4252
$ws.Range("N64").Value = -10493.5

$ws.Range("H67").Value = 6699
$ws.Range("I67").Value = 4500
$ws.Range("J67").Value = 9997.5
$ws.Range("K67").Value = 4500
$ws.Range("L67").Value = 9997.5
$ws.Range("M67").Value = -3642
$ws.Range("N67").Value = -11713.5

$ws.Range("H98").Value = 1409.5
$ws.Range("I98").Value = 1478.7142
$ws.Range("K98").Value = 1478.7142
$ws.Range("M98").Value = 19.28580000000011

$ws.Range("H100").Value = 3921
$ws.Range("I100").Value = 1317.8334
$ws.Range("K100").Value = 1317.8334
$ws.Range("M100").Value = -776.8334

$ws.Range("H103").Value = 669.875
$ws.Range("J103").Value = 871.5
$ws.Range("L103").Value = 2614.5
$ws.Range("N103").Value = -3786.5

$ws.Range("H106").Value = 5774.5713
$ws.Range("I106").Value = 5774.5713
$ws.Range("K106").Value = 5774.5713
$ws.Range("M106").Value = -5143.5713

$ws.Range("H113").Value = 4317.273
$ws.Range("I113").Value = 4186.25
$ws.Range("K113").Value = 4186.25
$ws.Range("M113").Value = -932.25

$ws.Range("H122").Value = 1409.5
$ws.Range("I122").Value = 1478.7142
$ws.Range("K122").Value = 4436.142599999999
$ws.Range("M122").Value = -1986.142599999999

$ws.Range("H132").Value = 36297.168
$ws.Range("I132").Value = 53458.25
$ws.Range("J132").Value = 1975
$ws.Range("K132").Value = 160374.75
$ws.Range("L132").Value = 5925
$ws.Range("M132").Value = -157844.75
$ws.Range("N132").Value = -10985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 142836.22
$ws.Range("I32").Value = 147591.12
$ws.Range("K32").Value = 147591.12
$ws.Range("M32").Value = -147304.12

$ws.Range("H45").Value = 3400

$ws.Range("H61").Value = 2409.1333
$ws.Range("I61").Value = 2379.7693
$ws.Range("J61").Value = 2600
$ws.Range("K61").Value = 2379.7693
$ws.Range("L61").Value = 2600
$ws.Range("M61").Value = -2167.7693
$ws.Range("N61").Value = -3024

$ws.Range("H110").Value = 1074
$ws.Range("I110").Value = 994.6667
$ws.Range("J110").Value = 1312
$ws.Range("K110").Value = 994.6667
$ws.Range("L110").Value = 1312
$ws.Range("M110").Value = 1050.3333
$ws.Range("N110").Value = -5402

$ws.Range("H132").Value = 8999.789000000001
$ws.Range("I132").Value = 6999.8
$ws.Range("J132").Value = 9714.071
$ws.Range("K132").Value = 20999.4
$ws.Range("L132").Value = 29142.213
$ws.Range("M132").Value = -18469.4
$ws.Range("N132").Value = -34202.213

$ws.Range("H136").Value = 2409.1333
$ws.Range("I136").Value = 2379.7693
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 7139.3079
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -4589.3079
$ws.Range("N136").Value = -12900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4281.591
$ws.Range("J94").Value = 4156.75
$ws.Range("L94").Value = 4156.75
$ws.Range("N94").Value = -5058.75

$ws.Range("H107").Value = 1264.7
$ws.Range("I107").Value = 1264.7
$ws.Range("K107").Value = 1264.7
$ws.Range("M107").Value = 655.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 187694.83
$ws.Range("I16").Value = 25233.8
$ws.Range("J16").Value = 1000000
$ws.Range("K16").Value = 25233.8
$ws.Range("L16").Value = 1000000
$ws.Range("M16").Value = -24946.8
$ws.Range("N16").Value = -1000574

$ws.Range("H99").Value = 17834.23
$ws.Range("I99").Value = 19012.084
$ws.Range("J99").Value = 3700
$ws.Range("K99").Value = 19012.084
$ws.Range("L99").Value = 3700
$ws.Range("M99").Value = -17514.084
$ws.Range("N99").Value = -6696

$ws.Range("H113").Value = 187694.83
$ws.Range("I113").Value = 25233.8
$ws.Range("J113").Value = 1000000
$ws.Range("K113").Value = 25233.8
$ws.Range("L113").Value = 1000000
$ws.Range("M113").Value = -23063.8
$ws.Range("N113").Value = -1004340

$ws.Range("H126").Value = 17834.23
$ws.Range("I126").Value = 19012.084
$ws.Range("J126").Value = 3700
$ws.Range("K126").Value = 57036.25199999999
$ws.Range("L126").Value = 11100
$ws.Range("M126").Value = -54566.25199999999
$ws.Range("N126").Value = -16040

$ws.Range("H134").Value = 3340.889
$ws.Range("I134").Value = 2814.8
$ws.Range("K134").Value = 8444.400000000001
$ws.Range("M134").Value = -5909.400000000001

$ws.Range("H141").Value = 47631.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 47631.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 47631.5
$ws.Range("N141").Value = -57991.5
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2923.75
$ws.Range("I6").Value = 377.75
$ws.Range("J6").Value = 5469.75
$ws.Range("K6").Value = 1133.25
$ws.Range("L6").Value = 16409.25
$ws.Range("M6").Value = -1020.25
$ws.Range("N6").Value = -16635.25

$ws.Range("H12").Value = 381.5357
$ws.Range("J12").Value = 201.34782
$ws.Range("L12").Value = 604.0434600000001
$ws.Range("N12").Value = -950.0434600000001

$ws.Range("H34").Value = 1522.1666
$ws.Range("J34").Value = 1992.8889
$ws.Range("L34").Value = 5978.6667
$ws.Range("N34").Value = -6146.6667

$ws.Range("H39").Value = 160891.6
$ws.Range("J39").Value = 79644.336
$ws.Range("L39").Value = 238933.008
$ws.Range("N39").Value = -239521.008

$ws.Range("H55").Value = 53127896
$ws.Range("J55").Value = 1115419
$ws.Range("L55").Value = 3346257
$ws.Range("N55").Value = -3346611

$ws.Range("H80").Value = 13849.9
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 13849.9
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 41549.7
$ws.Range("N80").Value = -43421.7
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 13849.9
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 13849.9
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 124649.1
$ws.Range("N83").Value = -134009.1
$ws.Range("M83").ClearContents()

$ws.Range("H113").Value = 25767.75
$ws.Range("I113").Value = 471.42856
$ws.Range("J113").Value = 36183.883
$ws.Range("K113").Value = 1414.28568
$ws.Range("L113").Value = 108551.649
$ws.Range("M113").Value = 755.71432
$ws.Range("N113").Value = -112891.649

$ws.Range("H122").Value = 1076420.4
$ws.Range("J122").Value = 1294.2307
$ws.Range("L122").Value = 11648.0763
$ws.Range("N122").Value = -16548.0763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 35807.395
$ws.Range("I2").Value = 58876.883
$ws.Range("J2").Value = 154.54546
$ws.Range("K2").Value = 58876.883
$ws.Range("L2").Value = 154.54546
$ws.Range("M2").Value = -58763.883
$ws.Range("N2").Value = -380.54546

$ws.Range("H122").Value = 1979.65
$ws.Range("J122").Value = 2999.5
$ws.Range("L122").Value = 8998.5
$ws.Range("N122").Value = -13898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1658.7916
$ws.Range("I16").Value = 1714.0454
$ws.Range("J16").Value = 1051
$ws.Range("K16").Value = 1714.0454
$ws.Range("L16").Value = 1051
$ws.Range("M16").Value = -1544.0454
$ws.Range("N16").Value = -1391

$ws.Range("H46").Value = 4297.905
$ws.Range("I46").Value = 1453.8
$ws.Range("K46").Value = 1453.8
$ws.Range("M46").Value = -1265.8

$ws.Range("H55").Value = 1440.5938
$ws.Range("I55").Value = 1365.8334
$ws.Range("J55").Value = 1485.45
$ws.Range("K55").Value = 1365.8334
$ws.Range("L55").Value = 1485.45
$ws.Range("M55").Value = -1192.8334
$ws.Range("N55").Value = -1831.45

$ws.Range("H61").Value = 9097.951999999999
$ws.Range("I61").Value = 9523.625
$ws.Range("J61").Value = 7735.8
$ws.Range("K61").Value = 9523.625
$ws.Range("L61").Value = 7735.8
$ws.Range("M61").Value = -9321.625
$ws.Range("N61").Value = -8139.8

$ws.Range("H113").Value = 9097.951999999999
$ws.Range("I113").Value = 9523.625
$ws.Range("J113").Value = 7735.8
$ws.Range("K113").Value = 9523.625
$ws.Range("L113").Value = 7735.8
$ws.Range("M113").Value = -7353.625
$ws.Range("N113").Value = -12075.8

$ws.Range("H132").Value = 2298.3215
$ws.Range("I132").Value = 1542.25
$ws.Range("J132").Value = 4188.5
$ws.Range("K132").Value = 4626.75
$ws.Range("L132").Value = 12565.5
$ws.Range("M132").Value = -2096.75
$ws.Range("N132").Value = -17625.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766

$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H107").Value = 1260.9556
$ws.Range("I107").Value = 937.62964
$ws.Range("J107").Value = 1745.9445
$ws.Range("K107").Value = 2812.88892
$ws.Range("L107").Value = 5237.833500000001
$ws.Range("M107").Value = -892.8889199999999
$ws.Range("N107").Value = -9077.833500000001

$ws.Range("H132").Value = 2466.4
$ws.Range("I132").Value = 2114.25
$ws.Range("J132").Value = 3875
$ws.Range("K132").Value = 6342.75
$ws.Range("L132").Value = 11625
$ws.Range("M132").Value = -3812.75
$ws.Range("N132").Value = -16685
